$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Handback generation: the a178a0eb-... file now has a target/handback file
# and an updated status + handback datetime, for both the zh-cn and de-de
# locale sheets. The e928ee0a-... rows are untouched.
# ---------------------------------------------------------------------------

$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("I2").Value = "a178a0eb-3d00-47e4-97bb-f2731cd9d2e6.md"
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d948bda97a9a363ee2631b84be63b875cf5a0d4b/e2e/a178a0eb-3d00-47e4-97bb-f2731cd9d2e6.md", "", "", "a178a0eb-3d00-47e4-97bb-f2731cd9d2e6.md")
$zh.Range("I2").Style = "Hyperlink"
$zh.Range("J2").Value = "a178a0eb-3d00-47e4-97bb-f2731cd9d2e6.6e9a706d6d15c846e6126abbd07546dbd0ac7efd.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-19 06:39:10"

$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("I2").Value = "a178a0eb-3d00-47e4-97bb-f2731cd9d2e6.md"
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d948bda97a9a363ee2631b84be63b875cf5a0d4b/e2e/a178a0eb-3d00-47e4-97bb-f2731cd9d2e6.md", "", "", "a178a0eb-3d00-47e4-97bb-f2731cd9d2e6.md")
$de.Range("I2").Style = "Hyperlink"
$de.Range("J2").Value = "a178a0eb-3d00-47e4-97bb-f2731cd9d2e6.6e9a706d6d15c846e6126abbd07546dbd0ac7efd.de-de.xlf"
$de.Range("K2").Value = "2016-08-19 06:39:17"

# ---------------------------------------------------------------------------
# Column widths widen to fit the new/longer content (Status column, and the
# newly-populated Latest Target File / Latest Handback File columns).
# ---------------------------------------------------------------------------

$ov = $wb.Worksheets.Item("Overview")
$ov.Columns.Item(5).ColumnWidth = 29.166666666666668
$ov.Columns.Item(6).ColumnWidth = 29.166666666666668

$zh.Columns.Item(3).ColumnWidth = 29.166666666666668
$zh.Columns.Item(9).ColumnWidth = 39.166666666666664
$zh.Columns.Item(10).ColumnWidth = 39.166666666666664

$de.Columns.Item(3).ColumnWidth = 29.166666666666668
$de.Columns.Item(9).ColumnWidth = 39.166666666666664
$de.Columns.Item(10).ColumnWidth = 39.166666666666664
